$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 5760, 45995),
    @(3, 5710, 45995.01041666666),
    @(4, 5690, 45995.02083333334),
    @(5, 5650, 45995.03125),
    @(6, 5600, 45995.04166666666),
    @(7, 5570, 45995.08333333334),
    @(8, 5580, 45995.11458333334),
    @(9, 5590, 45995.125),
    @(10, 5610, 45995.13541666666),
    @(11, 5640, 45995.14583333334),
    @(12, 5680, 45995.15625),
    @(13, 5740, 45995.16666666666),
    @(14, 5810, 45995.17708333334),
    @(15, 5920, 45995.1875),
    @(16, 6060, 45995.19791666666),
    @(17, 6230, 45995.20833333334),
    @(18, 6400, 45995.21875),
    @(19, 6590, 45995.22916666666),
    @(20, 6780, 45995.23958333334),
    @(21, 7030, 45995.25),
    @(22, 7220, 45995.26041666666),
    @(23, 7390, 45995.27083333334),
    @(24, 7550, 45995.28125),
    @(25, 7670, 45995.29166666666),
    @(26, 7770, 45995.30208333334),
    @(27, 7860, 45995.3125),
    @(28, 7900, 45995.32291666666),
    @(29, 7880, 45995.35416666666),
    @(30, 7820, 45995.36458333334),
    @(31, 7750, 45995.375),
    @(32, 7680, 45995.38541666666),
    @(33, 7610, 45995.39583333334),
    @(34, 7540, 45995.40625),
    @(35, 7460, 45995.41666666666),
    @(36, 7400, 45995.42708333334),
    @(37, 7360, 45995.4375),
    @(38, 7320, 45995.44791666666),
    @(39, 7290, 45995.45833333334),
    @(40, 7260, 45995.46875),
    @(41, 7250, 45995.47916666666),
    @(42, 7240, 45995.48958333334),
    @(43, 7250, 45995.5),
    @(44, 7270, 45995.51041666666),
    @(45, 7290, 45995.52083333334),
    @(46, 7300, 45995.53125),
    @(47, 7330, 45995.54166666666),
    @(48, 7370, 45995.55208333334),
    @(49, 7400, 45995.5625),
    @(50, 7430, 45995.57291666666),
    @(51, 7460, 45995.58333333334),
    @(52, 7500, 45995.59375),
    @(53, 7540, 45995.60416666666),
    @(54, 7590, 45995.61458333334),
    @(55, 7660, 45995.625),
    @(56, 7750, 45995.63541666666),
    @(57, 7850, 45995.64583333334),
    @(58, 7950, 45995.65625),
    @(59, 8080, 45995.66666666666),
    @(60, 8160, 45995.67708333334),
    @(61, 8200, 45995.6875),
    @(62, 8180, 45995.70833333334),
    @(63, 8140, 45995.71875),
    @(64, 8100, 45995.72916666666),
    @(65, 8070, 45995.73958333334),
    @(66, 8030, 45995.75),
    @(67, 8010, 45995.76041666666),
    @(68, 8000, 45995.77083333334),
    @(69, 7970, 45995.78125),
    @(70, 7920, 45995.79166666666),
    @(71, 7870, 45995.80208333334),
    @(72, 7790, 45995.8125),
    @(73, 7710, 45995.82291666666),
    @(74, 7600, 45995.83333333334),
    @(75, 7470, 45995.84375),
    @(76, 7370, 45995.85416666666),
    @(77, 7230, 45995.86458333334),
    @(78, 7050, 45995.875),
    @(79, 6920, 45995.88541666666),
    @(80, 6770, 45995.89583333334),
    @(81, 6630, 45995.90625),
    @(82, 6480, 45995.91666666666),
    @(83, 6350, 45995.92708333334),
    @(84, 6230, 45995.9375),
    @(85, 6120, 45995.94791666666),
    @(86, 5900, 45995.95833333334),
    @(87, 5830, 45995.96875),
    @(88, 5800, 45995.97916666666),
    @(89, 5750, 45995.98958333334)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
}

$ws.Range("A90:B91").EntireRow.Delete()
